$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-09-14"

$ws.Range("B1").Value = "September 2021 (through September 14)"

$ws.Range("K2").Value = 2
$ws.Range("BD3").Value = 1
$ws.Range("AU4").Value = 2
$ws.Range("K5").Value = 3
$ws.Range("AC8").Value = 1
$ws.Range("B13").Value = 4
$ws.Range("BD14").Value = 1
$ws.Range("K17").Value = 2
$ws.Range("H18").Value = 2
$ws.Range("AL18").Value = 3
$ws.Range("B19").Value = 3
$ws.Range("B27").Value = 2
$ws.Range("B30").Value = 1
$ws.Range("AL39").Value = 1
$ws.Range("AL41").Value = 1
$ws.Range("K43").Value = 2
$ws.Range("K46").Value = 3
$ws.Range("K55").Value = 3
$ws.Range("BD78").Value = 1
